# Revert "feat: rename acknowledgement of service (CMC-1271)"
# Change the heading text "Acknowledgement of claim" back to
# "Acknowledgement of service".

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "claim", $true, $true, $false, $false, $false,
    $true, 1, $false, "service", 2
)
